$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values (fix for int/float bug reading data from ByBit in find_crossing())
$ws.Range("A2").Value = 11
$ws.Range("B2").Value = "BTCUSD"

# Force text interpretation for the date-like and percent-like strings so Excel
# doesn't auto-convert them into a date serial / percentage number.
$ws.Range("C2").Value = "'2021-11-01"
$ws.Range("C2").ClearFormats()

$ws.Range("G2").Value = 1
$ws.Range("H2").Value = 0.6666666666666666
$ws.Range("K2").Value = $true
$ws.Range("L2").Value = 4
$ws.Range("M2").Value = 11
$ws.Range("N2").Value = 15

$ws.Range("O2").Value = "'26.7%"
$ws.Range("O2").ClearFormats()

$ws.Range("P2").Value = -10
$ws.Range("Q2").Value = 0
$ws.Range("R2").Value = 400
$ws.Range("S2").Value = -733.333333333333
$ws.Range("T2").Value = 185
$ws.Range("U2").Value = -518.333333333333

# Row 3 (Test #2) no longer exists in the backtest results - remove it entirely.
$ws.Rows("3:3").Delete()
